# Fruta / hortaliza, semanal
#
# Inserts a new daily price record (2021-12-20, "Castle Brite" / "Segunda")
# as row 4 of the data table, pushing the existing rows 4-30 down to 5-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 4 (shifts rows 4:30 -> 5:31).
$ws.Rows("4:4").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Macroferia Regional de Talca"
$ws.Range("C4").Value = "Maule"
$ws.Range("D4").Value = 44550
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103003
$ws.Range("J4").Value = "Damasco"
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "`$/caja 15 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 15
